# Implementation of ID tags in progress
#
# Adds seven new "ID tag" columns (D:J) to the LENGTHS sheet, with a header
# row of labels and single-letter tag values (A..L) down through row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LENGTHS")

# New header row (D1:J1)
$headers = @("TagID", "FloorID", "ZoneID", "LocationID", "MemberTypeID", "RebarTypeID", "SpecificTagID")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, 4 + $i).Value = $headers[$i]
}

# New column D values (D2:D13) - single letter tags, one per data row
$tags = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L")
for ($r = 0; $r -lt $tags.Count; $r++) {
    $ws.Cells.Item(2 + $r, 4).Value = $tags[$r]
}

# AutoFit the newly populated columns so their widths reflect the content,
# matching what Excel does automatically after data entry.
for ($c = 4; $c -le 10; $c++) {
    $ws.Columns.Item($c).EntireColumn.AutoFit() | Out-Null
}

# Move / update the active selection to the cell below the new data.
$ws.Range("D14").Select() | Out-Null
